$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update mac_address column (C) for existing rows 2-21 (Machine 1-20) with new MAC formatting
$ws.Cells.Item(2, 3).Value = "8C-16-45-5A-5D-0D"
$ws.Cells.Item(3, 3).Value = "8C-16-45-88-E1-0D"
$ws.Cells.Item(4, 3).Value = "00-FF-D3-E3-9A-27"
$ws.Cells.Item(5, 3).Value = "8C-16-45-5A-62-41"
$ws.Cells.Item(6, 3).Value = "E8-6A-64-1D-75-E4"
$ws.Cells.Item(7, 3).Value = "8C-16-45-FA-94-B7"
$ws.Cells.Item(8, 3).Value = "8C-16-45-1A-0F-62"
$ws.Cells.Item(9, 3).Value = "E8-6A-64-1C-52-6E"
$ws.Cells.Item(10, 3).Value = "48-51-B7-10-35-A6"
$ws.Cells.Item(11, 3).Value = "8C-16-45-38-F3-F3"
$ws.Cells.Item(12, 3).Value = "D4-3D-7E-58-CC-45"
$ws.Cells.Item(13, 3).Value = "8C-16-45-5A-5D-96"
$ws.Cells.Item(14, 3).Value = "8C-16-45-5A-5D-8E"
$ws.Cells.Item(15, 3).Value = "8C-16-45-33-A5-5F"
$ws.Cells.Item(16, 3).Value = "3C-95-09-F9-EA-DF"
$ws.Cells.Item(17, 3).Value = "8C-16-45-88-E7-0B"
$ws.Cells.Item(18, 3).Value = "B4-69-21-5A-DB-C4"
$ws.Cells.Item(19, 3).Value = "E8-6A-64-1D-48-B7"
$ws.Cells.Item(20, 3).Value = "8C-16-45-59-69-09 "
$ws.Cells.Item(21, 3).Value = "98-E7-F4-30-16-5A "

# Append new machine rows 22-30 (Machine 21-29)
$ws.Cells.Item(22, 1).Value = 10021
$ws.Cells.Item(22, 2).Value = "Machine 21"
$ws.Cells.Item(22, 3).Value = "38-BA-F8-53-C7-8F"
$ws.Cells.Item(22, 4).Value = "FB5962911653"
$ws.Cells.Item(22, 5).Value = "192.168.0.874"
$ws.Cells.Item(22, 6).Value = 1001
$ws.Cells.Item(22, 7).Value = "eng"
$ws.Cells.Item(22, 8).Value = $true
$ws.Cells.Item(22, 9).Value = "superadmin"
$ws.Cells.Item(22, 10).Value = "now()"
$ws.Cells.Item(22, 11).Value = "now()"

$ws.Cells.Item(23, 1).Value = 10022
$ws.Cells.Item(23, 2).Value = "Machine 22"
$ws.Cells.Item(23, 3).Value = "E8-6A-64-1C-58-C2"
$ws.Cells.Item(23, 4).Value = "FB5962911654"
$ws.Cells.Item(23, 5).Value = "192.168.0.721"
$ws.Cells.Item(23, 6).Value = 1001
$ws.Cells.Item(23, 7).Value = "eng"
$ws.Cells.Item(23, 8).Value = $true
$ws.Cells.Item(23, 9).Value = "superadmin"
$ws.Cells.Item(23, 10).Value = "now()"
$ws.Cells.Item(23, 11).Value = "now()"

$ws.Cells.Item(24, 1).Value = 10023
$ws.Cells.Item(24, 2).Value = "Machine 23"
$ws.Cells.Item(24, 3).Value = "E4-A4-71-CE-BA-93"
$ws.Cells.Item(24, 4).Value = "FB5962911655"
$ws.Cells.Item(24, 5).Value = "192.168.0.841"
$ws.Cells.Item(24, 6).Value = 1001
$ws.Cells.Item(24, 7).Value = "eng"
$ws.Cells.Item(24, 8).Value = $true
$ws.Cells.Item(24, 9).Value = "superadmin"
$ws.Cells.Item(24, 10).Value = "now()"
$ws.Cells.Item(24, 11).Value = "now()"

$ws.Cells.Item(25, 1).Value = 10024
$ws.Cells.Item(25, 2).Value = "Machine 24"
$ws.Cells.Item(25, 3).Value = "54-E1-AD-EA-30-C9"
$ws.Cells.Item(25, 4).Value = "FB5962911656"
$ws.Cells.Item(25, 5).Value = "192.168.0.186"
$ws.Cells.Item(25, 6).Value = 1001
$ws.Cells.Item(25, 7).Value = "eng"
$ws.Cells.Item(25, 8).Value = $true
$ws.Cells.Item(25, 9).Value = "superadmin"
$ws.Cells.Item(25, 10).Value = "now()"
$ws.Cells.Item(25, 11).Value = "now()"

$ws.Cells.Item(26, 1).Value = 10025
$ws.Cells.Item(26, 2).Value = "Machine 25"
$ws.Cells.Item(26, 3).Value = "8C-16-45-65-DD-40"
$ws.Cells.Item(26, 4).Value = "FB5962911657"
$ws.Cells.Item(26, 5).Value = "192.168.0.627"
$ws.Cells.Item(26, 6).Value = 1001
$ws.Cells.Item(26, 7).Value = "eng"
$ws.Cells.Item(26, 8).Value = $true
$ws.Cells.Item(26, 9).Value = "superadmin"
$ws.Cells.Item(26, 10).Value = "now()"
$ws.Cells.Item(26, 11).Value = "now()"

$ws.Cells.Item(27, 1).Value = 10026
$ws.Cells.Item(27, 2).Value = "Machine 26"
$ws.Cells.Item(27, 3).Value = "58-20-B1-D6-C3-BE"
$ws.Cells.Item(27, 4).Value = "FB5962911658"
$ws.Cells.Item(27, 5).Value = "192.168.0.879"
$ws.Cells.Item(27, 6).Value = 1001
$ws.Cells.Item(27, 7).Value = "eng"
$ws.Cells.Item(27, 8).Value = $true
$ws.Cells.Item(27, 9).Value = "superadmin"
$ws.Cells.Item(27, 10).Value = "now()"
$ws.Cells.Item(27, 11).Value = "now()"

$ws.Cells.Item(28, 1).Value = 10027
$ws.Cells.Item(28, 2).Value = "Machine 27"
$ws.Cells.Item(28, 3).Value = "8C-16-45-38-F0-25"
$ws.Cells.Item(28, 4).Value = "FB5962911659"
$ws.Cells.Item(28, 5).Value = "192.168.0.628"
$ws.Cells.Item(28, 6).Value = 1001
$ws.Cells.Item(28, 7).Value = "eng"
$ws.Cells.Item(28, 8).Value = $true
$ws.Cells.Item(28, 9).Value = "superadmin"
$ws.Cells.Item(28, 10).Value = "now()"
$ws.Cells.Item(28, 11).Value = "now()"

$ws.Cells.Item(29, 1).Value = 10028
$ws.Cells.Item(29, 2).Value = "Machine 28"
$ws.Cells.Item(29, 3).Value = "6C-88-14-AC-EF-55"
$ws.Cells.Item(29, 4).Value = "FB5962911661"
$ws.Cells.Item(29, 5).Value = "192.168.0.306"
$ws.Cells.Item(29, 6).Value = 1001
$ws.Cells.Item(29, 7).Value = "eng"
$ws.Cells.Item(29, 8).Value = $true
$ws.Cells.Item(29, 9).Value = "superadmin"
$ws.Cells.Item(29, 10).Value = "now()"
$ws.Cells.Item(29, 11).Value = "now()"

$ws.Cells.Item(30, 1).Value = 10029
$ws.Cells.Item(30, 2).Value = "Machine 29"
$ws.Cells.Item(30, 3).Value = "3C-6A-A7-C0-DF-27"
$ws.Cells.Item(30, 4).Value = "FB5962911662"
$ws.Cells.Item(30, 5).Value = "192.168.0.355"
$ws.Cells.Item(30, 6).Value = 1001
$ws.Cells.Item(30, 7).Value = "eng"
$ws.Cells.Item(30, 8).Value = $true
$ws.Cells.Item(30, 9).Value = "superadmin"
$ws.Cells.Item(30, 10).Value = "now()"
$ws.Cells.Item(30, 11).Value = "now()"

# Column C width adjustment (to fit new MAC address format)
$ws.Columns.Item(3).ColumnWidth = 16.16666666666667

# Update selection to reflect the blank area below the data (as in the saved file)
$ws.Range("A31:XFD1048576").Select()